$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "messageType" column is gone; remaining header row is re-packed and a new
# "commTime" column is appended at the end.
$ws.Cells.Item(1, 3).Value  = "sender"
$ws.Cells.Item(1, 4).Value  = "receiver"
$ws.Cells.Item(1, 5).Value  = "messageContent"
$ws.Cells.Item(1, 6).Value  = "messageId"
$ws.Cells.Item(1, 7).Value  = "messageSize"
$ws.Cells.Item(1, 8).Value  = "commProcId"
$ws.Cells.Item(1, 9).Value  = "commStatus"
$ws.Cells.Item(1, 10).Value = "commTime"

# Select the whole header row, matching the saved selection state.
[void]$ws.Range("A1:XFD1").Select()

# Set up the page for printing (paper size 9 = A4, portrait).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
